# contingencies with rene fine
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Formatting: bold font, centered/top-aligned, thin box border around B1
$cellB1 = $ws.Range("B1")
$cellB1.Font.Bold = $true
$cellB1.HorizontalAlignment = -4108
$cellB1.VerticalAlignment = -4160
$cellB1.Borders.LineStyle = 1
$cellB1.Borders.Weight = 2

# Apply the exact same style to A2 by copying formats from B1,
# so both cells share a single style entry instead of generating duplicates.
$cellB1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
